$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "23.527.72"
Set-TextValue $ws.Range("E2") "  +2.29%  "
Set-TextValue $ws.Range("D3") "1.633.52"
Set-TextValue $ws.Range("E3") "  +3.15%  "
Set-TextValue $ws.Range("D4") "0.9956"
Set-TextValue $ws.Range("E4") "  -0.79%  "
Set-TextValue $ws.Range("D5") "307.98"
Set-TextValue $ws.Range("E5") "  +2.85%  "
Set-TextValue $ws.Range("D6") "0.9967"
Set-TextValue $ws.Range("E6") "  -0.65%  "
Set-TextValue $ws.Range("D7") "0.3786"
Set-TextValue $ws.Range("E7") "  +0.96%  "
Set-TextValue $ws.Range("D8") "53.19"
Set-TextValue $ws.Range("E8") "  +5.37%  "
Set-TextValue $ws.Range("D9") "0.3668"
Set-TextValue $ws.Range("E9") "  +3.01%  "
Set-TextValue $ws.Range("D10") "1.292"
Set-TextValue $ws.Range("E10") "  +6.70%  "
Set-TextValue $ws.Range("D11") "0.08205"
Set-TextValue $ws.Range("E11") "  +3.01%  "
Set-TextValue $ws.Range("D12") "0.9959"
Set-TextValue $ws.Range("E12") "  -0.79%  "
Set-TextValue $ws.Range("D13") "23.36"
Set-TextValue $ws.Range("E13") "  +7.65%  "
Set-TextValue $ws.Range("D14") "6.689"
Set-TextValue $ws.Range("E14") "  +3.91%  "
Set-TextValue $ws.Range("D15") "0.00001275"
Set-TextValue $ws.Range("E15") "  +4.93%  "
Set-TextValue $ws.Range("D16") "7.479"
Set-TextValue $ws.Range("E16") "  +2.97%  "
Set-TextValue $ws.Range("D17") "1.630.04"
Set-TextValue $ws.Range("E17") "  +2.75%  "
Set-TextValue $ws.Range("D18") "94.98"
Set-TextValue $ws.Range("E18") "  +3.39%  "
Set-TextValue $ws.Range("D19") "0.06946"
Set-TextValue $ws.Range("E19") "  +2.95%  "
Set-TextValue $ws.Range("D20") "18.46"
Set-TextValue $ws.Range("E20") "  +4.33%  "
Set-TextValue $ws.Range("D21") "6.608"
Set-TextValue $ws.Range("E21") "  +3.71%  "
Set-TextValue $ws.Range("D22") "0.9972"
Set-TextValue $ws.Range("E22") "  -0.57%  "
Set-TextValue $ws.Range("D23") "13.03"
Set-TextValue $ws.Range("E23") "  +2.68%  "
Set-TextValue $ws.Range("D24") "23.530.12"
Set-TextValue $ws.Range("E24") "  +2.31%  "
Set-TextValue $ws.Range("D25") "3.139"
Set-TextValue $ws.Range("E25") "  +13.68%  "
Set-TextValue $ws.Range("D26") "2.439"
Set-TextValue $ws.Range("E26") "  +2.88%  "
Set-TextValue $ws.Range("D27") "21.45"
Set-TextValue $ws.Range("E27") "  +4.17%  "
Set-TextValue $ws.Range("D28") "151.36"
Set-TextValue $ws.Range("E28") "  +3.29%  "
Set-TextValue $ws.Range("D29") "5.317"
Set-TextValue $ws.Range("E29") "  +2.18%  "
Set-TextValue $ws.Range("D30") "136.55"
Set-TextValue $ws.Range("E30") "  +3.76%  "
Set-TextValue $ws.Range("D31") "2.424"
Set-TextValue $ws.Range("E31") "  +4.95%  "
Set-TextValue $ws.Range("D32") "6.929"
Set-TextValue $ws.Range("E32") "  +7.19%  "
Set-TextValue $ws.Range("D33") "1.806.36"
Set-TextValue $ws.Range("E33") "  +2.19%  "
Set-TextValue $ws.Range("D34") "0.9827"
Set-TextValue $ws.Range("E34") "  +6.00%  "
Set-TextValue $ws.Range("D35") "0.02812"
Set-TextValue $ws.Range("E35") "  +5.93%  "
Set-TextValue $ws.Range("D36") "10.49"
Set-TextValue $ws.Range("E36") "  +6.54%  "
Set-TextValue $ws.Range("D37") "0.07491"
Set-TextValue $ws.Range("E37") "  +2.53%  "
Set-TextValue $ws.Range("D38") "6.292"
Set-TextValue $ws.Range("E38") "  +5.34%  "
Set-TextValue $ws.Range("D39") "0.2541"
Set-TextValue $ws.Range("E39") "  +2.54%  "
Set-TextValue $ws.Range("D40") "0.08856"
Set-TextValue $ws.Range("E40") "  +1.45%  "
Set-TextValue $ws.Range("D41") "1.410"
Set-TextValue $ws.Range("E41") "  +5.98%  "
Set-TextValue $ws.Range("D42") "0.7192"
Set-TextValue $ws.Range("E42") "  +5.57%  "
Set-TextValue $ws.Range("D43") "12.84"
Set-TextValue $ws.Range("E43") "  +9.79%  "
Set-TextValue $ws.Range("D44") "16.31"
Set-TextValue $ws.Range("E44") "  +10.90%  "
Set-TextValue $ws.Range("D45") "0.6642"
Set-TextValue $ws.Range("E45") "  +5.01%  "
Set-TextValue $ws.Range("D46") "2.371"
Set-TextValue $ws.Range("E46") "  +6.41%  "
Set-TextValue $ws.Range("D47") "4.041"
Set-TextValue $ws.Range("E47") "  +1.82%  "
Set-TextValue $ws.Range("D48") "0.9964"
Set-TextValue $ws.Range("E48") "  -0.56%  "
Set-TextValue $ws.Range("D49") "0.08051"
Set-TextValue $ws.Range("E49") "  +2.37%  "
Set-TextValue $ws.Range("D50") "132.93"
Set-TextValue $ws.Range("E50") "  +1.78%  "
Set-TextValue $ws.Range("E51") "  +3.61%  "
